$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "256.63"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "0.01%"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "15"
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "26.59"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "-1.87%"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "15"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "4.639"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "0.88%"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "15"
$ws.Cells.Item(4, 7).Style = "Normal"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.05928"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "0.44%"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "15"
$ws.Cells.Item(5, 7).Style = "Normal"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "6.602"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "-0.64%"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "15"
$ws.Cells.Item(6, 7).Style = "Normal"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.8553"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "-1.48%"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "15"
$ws.Cells.Item(7, 7).Style = "Normal"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.9134"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "-2.58%"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "15"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "One"
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.0006024"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "0.33%"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "15"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 2).Style = "Normal"
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.1375"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "-1.85%"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "15"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.04464"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "15.50%"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "15"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12, 2).Style = "Normal"
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.06992"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "-1.15%"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "15"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 2).Style = "Normal"
$ws.Cells.Item(13, 3).NumberFormat = "@"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.03014"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "-5.82%"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "15"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 2).Style = "Normal"
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09109"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "-1.61%"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "15"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "BitForexToken"
$ws.Cells.Item(15, 2).Style = "Normal"
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.001525"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "-1.88%"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "15"
$ws.Cells.Item(15, 7).Style = "Normal"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.006043"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "-1.02%"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "15"
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.468"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "-1.25%"
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "15"
$ws.Cells.Item(17, 7).Style = "Normal"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.132"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "-1.77%"
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "15"
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "-2.36%"
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "15"
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.3075"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "0.12%"
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "15"
$ws.Cells.Item(20, 7).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "0.96%"
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "15"
$ws.Cells.Item(21, 7).Style = "Normal"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.850"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "0.12%"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "15"
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04196"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "-0.29%"
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "15"
$ws.Cells.Item(23, 7).Style = "Normal"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.001215"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-0.04%"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "15"
$ws.Cells.Item(24, 7).Style = "Normal"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.004473"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "4.55%"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "15"
$ws.Cells.Item(25, 7).Style = "Normal"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.0001199"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "0.03%"
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "15"
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0001713"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "13.76%"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "15"
$ws.Cells.Item(27, 7).Style = "Normal"
$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "15"
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "15"
$ws.Cells.Item(29, 7).Style = "Normal"
$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "15"
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "15"
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "15"
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "15"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "15"
$ws.Cells.Item(34, 7).Style = "Normal"
$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "15"
$ws.Cells.Item(35, 7).Style = "Normal"
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "15"
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "15"
$ws.Cells.Item(37, 7).Style = "Normal"
$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "15"
$ws.Cells.Item(38, 7).Style = "Normal"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "15"
$ws.Cells.Item(39, 7).Style = "Normal"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.03795"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "-0.68%"
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "15"
$ws.Cells.Item(40, 7).Style = "Normal"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.006213"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "57.97%"
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "15"
$ws.Cells.Item(41, 7).Style = "Normal"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1098"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "-0.13%"
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "15"
$ws.Cells.Item(42, 7).Style = "Normal"
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = "LocalTraders"
$ws.Cells.Item(43, 2).Style = "Normal"
$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.01499"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "27.74%"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "15"
$ws.Cells.Item(43, 7).Style = "Normal"
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = "CEJI"
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.002145"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "-6.77%"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "15"
$ws.Cells.Item(44, 7).Style = "Normal"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.00005082"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "-6.65%"
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "15"
$ws.Cells.Item(45, 7).Style = "Normal"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00000000749"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "15"
$ws.Cells.Item(46, 7).Style = "Normal"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.04995"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "-35.74%"
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "15"
$ws.Cells.Item(47, 7).Style = "Normal"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "10,483.54%"
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "15"
$ws.Cells.Item(48, 7).Style = "Normal"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.00002098"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "15"
$ws.Cells.Item(49, 7).Style = "Normal"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0001998"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "15"
$ws.Cells.Item(50, 7).Style = "Normal"
$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "15"
$ws.Cells.Item(51, 7).Style = "Normal"
